# edit.ps1
# Applies the data edits described by the commit diff to before.xlsx:
#  - Sheet "Inventory": updates rows 2-5, appends rows 6-11 (dimension A1:O5 -> A1:O11)
#  - Sheet "Sales": appends rows 16-17 (dimension A1:F15 -> A1:F17)

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Inventory")
$ws2 = $wb.Worksheets.Item("Sales")

# Helper: force a numeric-looking value to be written as TEXT (matching the
# source workbook's convention of storing many numeric-looking fields as
# strings). A direct `$cell.Value = "5.56"` assignment auto-coerces to a
# number (standard Excel behavior), so instead we stage the text in a
# scratch cell that is explicitly number-formatted as Text ("@"), copy it,
# and paste-special (values only) into the destination -- this keeps the
# destination's existing formatting/style intact while forcing a text type.
# The scratch row is removed afterwards so it doesn't affect the sheet's
# used range / dimension.
function Set-TextValue($Worksheet, $Row, $Col, $Text) {
    $scratch = $Worksheet.Cells.Item(900, 1)
    $scratch.NumberFormat = "@"
    $scratch.Value = $Text
    $scratch.Copy()
    $Worksheet.Cells.Item($Row, $Col).PasteSpecial(-4163)
    $scratch.ClearContents()
    $scratch.NumberFormat = "General"
    $scratch.EntireRow.Delete()
}

# ===== Sheet1 (Inventory): update existing rows 2-5 =====
$ws1.Cells.Item(2,1).Value = "glop"
$ws1.Cells.Item(2,3).Value = "eth"
$ws1.Cells.Item(2,4).Value = 12
$ws1.Cells.Item(2,8).Value = 10.5
$ws1.Cells.Item(2,13).Value = 23

$ws1.Cells.Item(3,1).Value = "glplpp"
$ws1.Cells.Item(3,13).Value = 21

$ws1.Cells.Item(4,1).Value = "emp"
$ws1.Cells.Item(4,3).Value = "chn"
$ws1.Cells.Item(4,4).Value = 10
$ws1.Cells.Item(4,6).Value = 249
$ws1.Cells.Item(4,7).Value = 951235648952
$ws1.Cells.Item(4,8).Value = 12.35
Set-TextValue $ws1 4 9 "5.56"
Set-TextValue $ws1 4 10 "17.91"
Set-TextValue $ws1 4 11 "2.69"
Set-TextValue $ws1 4 12 "20.59"
$ws1.Cells.Item(4,13).Value = 201
$ws1.Cells.Item(4,14).Value = 0
$ws1.Cells.Item(4,15).Value = 0

$ws1.Cells.Item(5,1).Value = "opads"
Set-TextValue $ws1 5 4 "10"
Set-TextValue $ws1 5 8 "12.35"
Set-TextValue $ws1 5 9 "5.4"
Set-TextValue $ws1 5 10 "17.40"
Set-TextValue $ws1 5 11 "2.61"
Set-TextValue $ws1 5 12 "20.01"
$ws1.Cells.Item(5,13).Value = 239

# ===== Sheet1 (Inventory): append new rows 6-11 =====
# -- row 6 --
$ws1.Cells.Item(6,1).Value = "lipbalm"
$ws1.Cells.Item(6,2).Value = "lipstick"
$ws1.Cells.Item(6,3).Value = "chn"
$ws1.Cells.Item(6,4).Value = 10
$ws1.Cells.Item(6,5).Value = "pcs"
$ws1.Cells.Item(6,6).Value = 34
$ws1.Cells.Item(6,7).Value = 951235648952
$ws1.Cells.Item(6,8).Value = 10.4
Set-TextValue $ws1 6 9 "4.68"
Set-TextValue $ws1 6 10 "15.08"
Set-TextValue $ws1 6 11 "2.26"
Set-TextValue $ws1 6 12 "17.34"
$ws1.Cells.Item(6,13).Value = 34
$ws1.Cells.Item(6,14).Value = 0
$ws1.Cells.Item(6,15).Value = 0

# -- row 7 --
$ws1.Cells.Item(7,1).Value = "asfafas"
$ws1.Cells.Item(7,2).Value = "stock"
$ws1.Cells.Item(7,3).Value = "chn"
$ws1.Cells.Item(7,4).Value = 10
$ws1.Cells.Item(7,5).Value = "pcs"
$ws1.Cells.Item(7,6).Value = 249
$ws1.Cells.Item(7,7).Value = 951235648952
$ws1.Cells.Item(7,8).Value = 12.35
Set-TextValue $ws1 7 9 "5.56"
Set-TextValue $ws1 7 10 "17.91"
Set-TextValue $ws1 7 11 "2.69"
Set-TextValue $ws1 7 12 "20.59"
$ws1.Cells.Item(7,13).Value = 247
$ws1.Cells.Item(7,14).Value = 0
$ws1.Cells.Item(7,15).Value = 0

# -- row 8 --
$ws1.Cells.Item(8,1).Value = "etett"
$ws1.Cells.Item(8,2).Value = "stock"
$ws1.Cells.Item(8,3).Value = "chn"
$ws1.Cells.Item(8,4).Value = 10
$ws1.Cells.Item(8,5).Value = "pcs"
$ws1.Cells.Item(8,6).Value = 249
$ws1.Cells.Item(8,7).Value = 951235648952
$ws1.Cells.Item(8,8).Value = 454
Set-TextValue $ws1 8 9 "204.30"
Set-TextValue $ws1 8 10 "658.30"
Set-TextValue $ws1 8 11 "98.74"
Set-TextValue $ws1 8 12 "757.04"
$ws1.Cells.Item(8,13).Value = 249
$ws1.Cells.Item(8,14).Value = 0
$ws1.Cells.Item(8,15).Value = 0

# -- row 9 --
$ws1.Cells.Item(9,1).Value = "empkmkm"
$ws1.Cells.Item(9,2).Value = "stock"
$ws1.Cells.Item(9,3).Value = "chn"
$ws1.Cells.Item(9,4).Value = 10
$ws1.Cells.Item(9,5).Value = "pcs"
$ws1.Cells.Item(9,6).Value = 249
$ws1.Cells.Item(9,7).Value = 951235648952
$ws1.Cells.Item(9,8).Value = 12.32
Set-TextValue $ws1 9 9 "5.54"
Set-TextValue $ws1 9 10 "17.86"
Set-TextValue $ws1 9 11 "2.68"
Set-TextValue $ws1 9 12 "20.54"
$ws1.Cells.Item(9,13).Value = 249
$ws1.Cells.Item(9,14).Value = 0
$ws1.Cells.Item(9,15).Value = 0

# -- row 10 --
$ws1.Cells.Item(10,1).Value = "emp"
$ws1.Cells.Item(10,2).Value = "stock"
$ws1.Cells.Item(10,3).Value = "chn"
$ws1.Cells.Item(10,4).Value = 10
$ws1.Cells.Item(10,5).Value = "pcs"
$ws1.Cells.Item(10,6).Value = 249
$ws1.Cells.Item(10,7).Value = 951235648952
$ws1.Cells.Item(10,8).Value = 12.35
Set-TextValue $ws1 10 9 "5.56"
Set-TextValue $ws1 10 10 "17.91"
Set-TextValue $ws1 10 11 "2.69"
Set-TextValue $ws1 10 12 "20.59"
$ws1.Cells.Item(10,13).Value = 249
$ws1.Cells.Item(10,14).Value = 0
$ws1.Cells.Item(10,15).Value = 0

# -- row 11 --
$ws1.Cells.Item(11,1).Value = "emp"
$ws1.Cells.Item(11,2).Value = "lipstick"
$ws1.Cells.Item(11,3).Value = "chn"
Set-TextValue $ws1 11 4 "10"
$ws1.Cells.Item(11,5).Value = "pcs"
Set-TextValue $ws1 11 6 "249"
Set-TextValue $ws1 11 7 "951235648952"
Set-TextValue $ws1 11 8 "12.35"
Set-TextValue $ws1 11 9 "5.56"
Set-TextValue $ws1 11 10 "17.91"
Set-TextValue $ws1 11 11 "2.69"
Set-TextValue $ws1 11 12 "20.59"
Set-TextValue $ws1 11 13 "249"
$ws1.Cells.Item(11,14).Value = 0
$ws1.Cells.Item(11,15).Value = 0

# ===== Sheet2 (Sales): append new rows 16-17 =====
# -- row 16 --
$ws2.Cells.Item(16,1).Value = "opads"
$ws2.Cells.Item(16,2).Value = 12.35
$ws2.Cells.Item(16,3).Value = 10
$ws2.Cells.Item(16,4).Value = 123.5
$ws2.Cells.Item(16,5).Value = "selome"
$ws2.Cells.Item(16,6).Value = "2023-01-19 15:56:08"

# -- row 17 --
$ws2.Cells.Item(17,1).Value = "asfafas"
$ws2.Cells.Item(17,2).Value = 12.35
$ws2.Cells.Item(17,3).Value = 2
$ws2.Cells.Item(17,4).Value = 24.7
$ws2.Cells.Item(17,5).Value = "selome"
$ws2.Cells.Item(17,6).Value = "2023-01-19 16:32:50"


Write-Output "Edit complete: Inventory rows 2-11, Sales rows 16-17 updated."
